$wb = $excel.ActiveWorkbook

# Insert the new "CreditCardDetails" worksheet directly before "AccountCreationData".
# Excel.Worksheets.Add(Before) mirrors the VBA signature Add(Before, After, Count, Type).
$accountCreationSheet = $wb.Worksheets.Item("AccountCreationData")
$ws = $wb.Worksheets.Add($accountCreationSheet)
$ws.Name = "CreditCardDetails"

# Header row
$ws.Range("A1").Value = "CCHolderName"
$ws.Range("B1").Value = "CCNumber"
$ws.Range("C1").Value = "CVC"
$ws.Range("D1").Value = "CCExpiryMonth"
$ws.Range("E1").Value = "CCExpiryYear"

# Data row
$ws.Range("A2").Value = "Test Automation"
$ws.Range("B2").Value = "4111 1111 1111 1111"
$ws.Range("C2").Value = 737
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 2030

# Formatting: bold header with a thin box border all around, same border on data row.
$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1

$dataRange = $ws.Range("A2:E2")
$dataRange.Borders.LineStyle = 1

# Column widths to roughly match the authored sheet (A, B wider for names/numbers; D, E for dates)
$ws.Columns.Item(1).ColumnWidth = 13.833333333333334
$ws.Columns.Item(2).ColumnWidth = 17.5
$ws.Columns.Item(4).ColumnWidth = 13.0
$ws.Columns.Item(5).ColumnWidth = 13.0

# Make the new sheet the active tab, with F17 selected (matches the authored selection).
$ws.Activate()
[void]$ws.Range("F17").Select()
